$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row at row 744, shifting existing rows (744:936) down to (745:937)
$ws.Rows.Item(744).Insert()

# Populate the newly inserted row with the manually-added review data
$ws.Range("A744").Value = -1
$ws.Range("B744").Value = "uber crooked"
$ws.Range("C744").NumberFormat = "@"
$ws.Range("C744").Value = "-0.87784"
$ws.Range("D744").Value = -1

# Update the view state (scroll position + active cell) to match the author's final selection
$ws.Range("E744").Select()
